# Scheduled runner update: refresh currentAveragePrice / Leve profit figures
# for a batch of leve rows across several sheets.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 498445.44
$ws.Range("I129").Value = 445.4
$ws.Range("J129").Value = 582852.25
$ws.Range("K129").Value = 1336.2
$ws.Range("L129").Value = 1748556.75
$ws.Range("M129").Value = 3663.8
$ws.Range("N129").Value = -1758556.75

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4023.746
$ws.Range("I32").Value = 3964.1018
$ws.Range("J32").Value = 4903.5
$ws.Range("K32").Value = 3964.1018
$ws.Range("L32").Value = 4903.5
$ws.Range("M32").Value = -3677.1018
$ws.Range("N32").Value = -5477.5

$ws.Range("H132").Value = 6716.087
$ws.Range("I132").Value = 9089.23
$ws.Range("J132").Value = 3631
$ws.Range("K132").Value = 27267.69
$ws.Range("L132").Value = 10893
$ws.Range("M132").Value = -24737.69
$ws.Range("N132").Value = -15953

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H56").Value = 17000
$ws.Range("J56").Value = 17000
$ws.Range("L56").Value = 17000
$ws.Range("N56").Value = -18478

$ws.Range("H99").Value = 1155
$ws.Range("I99").Value = 1003.8
$ws.Range("K99").Value = 1003.8
$ws.Range("M99").Value = 494.2

$ws.Range("H105").Value = 2402.7273
$ws.Range("I105").Value = 1763.625
$ws.Range("J105").Value = 4107
$ws.Range("K105").Value = 1763.625
$ws.Range("L105").Value = 4107
$ws.Range("M105").Value = -16.625
$ws.Range("N105").Value = -7601

$ws.Range("H134").Value = 45414.914
$ws.Range("I134").Value = 64027.688
$ws.Range("J134").Value = 2871.4285
$ws.Range("K134").Value = 192083.064
$ws.Range("L134").Value = 8614.2855
$ws.Range("M134").Value = -189548.064
$ws.Range("N134").Value = -13684.2855

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 978.1818
$ws.Range("I16").Value = 894.2857
$ws.Range("J16").Value = 1125
$ws.Range("K16").Value = 894.2857
$ws.Range("L16").Value = 1125
$ws.Range("M16").Value = -607.2857
$ws.Range("N16").Value = -1699

$ws.Range("H31").Value = 6669637.5
$ws.Range("I31").Value = 3111.7144
$ws.Range("J31").Value = 100001000
$ws.Range("K31").Value = 3111.7144
$ws.Range("L31").Value = 100001000
$ws.Range("M31").Value = -2816.7144
$ws.Range("N31").Value = -100001590

$ws.Range("H34").Value = 6669637.5
$ws.Range("I34").Value = 3111.7144
$ws.Range("J34").Value = 100001000
$ws.Range("K34").Value = 3111.7144
$ws.Range("L34").Value = 100001000
$ws.Range("M34").Value = -2909.7144
$ws.Range("N34").Value = -100001404

$ws.Range("H58").Value = 16668037
$ws.Range("I58").Value = 1521
$ws.Range("J58").Value = 30304276
$ws.Range("K58").Value = 1521
$ws.Range("L58").Value = 30304276
$ws.Range("M58").Value = -1318
$ws.Range("N58").Value = -30304682

$ws.Range("H94").Value = 3220
$ws.Range("J94").Value = 3220
$ws.Range("L94").Value = 3220
$ws.Range("N94").Value = -4122

$ws.Range("H105").Value = 825.64703
$ws.Range("I105").Value = 681.6667
$ws.Range("J105").Value = 987.625
$ws.Range("K105").Value = 681.6667
$ws.Range("L105").Value = 987.625
$ws.Range("M105").Value = 1065.3333
$ws.Range("N105").Value = -4481.625

$ws.Range("H107").Value = 501
$ws.Range("I107").Value = 439.6
$ws.Range("J107").Value = 720.2857
$ws.Range("K107").Value = 439.6
$ws.Range("L107").Value = 720.2857
$ws.Range("M107").Value = 1480.4
$ws.Range("N107").Value = -4560.2857

$ws.Range("H113").Value = 978.1818
$ws.Range("I113").Value = 894.2857
$ws.Range("J113").Value = 1125
$ws.Range("K113").Value = 894.2857
$ws.Range("L113").Value = 1125
$ws.Range("M113").Value = 1275.7143
$ws.Range("N113").Value = -5465

$ws.Range("H134").Value = 1152.5
$ws.Range("I134").Value = 1109.0667
$ws.Range("J134").Value = 1282.8
$ws.Range("K134").Value = 3327.2001
$ws.Range("L134").Value = 3848.4
$ws.Range("M134").Value = -792.2001
$ws.Range("N134").Value = -8918.4

$ws.Range("H136").Value = 16668037
$ws.Range("I136").Value = 1521
$ws.Range("J136").Value = 30304276
$ws.Range("K136").Value = 4563
$ws.Range("L136").Value = 90912828
$ws.Range("M136").Value = -2013
$ws.Range("N136").Value = -90917928

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 788.4545000000001
$ws.Range("I5").Value = 191.6
$ws.Range("J5").Value = 1285.8334
$ws.Range("K5").Value = 574.8
$ws.Range("L5").Value = 3857.5002
$ws.Range("M5").Value = -462.8
$ws.Range("N5").Value = -4081.5002

$ws.Range("H118").Value = 4666.6665
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").Value = $null

$ws.Range("H122").Value = 1449.9524
$ws.Range("I122").Value = 1365.1538
$ws.Range("J122").Value = 1587.75
$ws.Range("K122").Value = 12286.3842
$ws.Range("L122").Value = 14289.75
$ws.Range("M122").Value = -9836.3842
$ws.Range("N122").Value = -19189.75

$ws.Range("H135").Value = 788.4545000000001
$ws.Range("I135").Value = 191.6
$ws.Range("J135").Value = 1285.8334
$ws.Range("K135").Value = 1724.4
$ws.Range("L135").Value = 11572.5006
$ws.Range("M135").Value = 810.6000000000001
$ws.Range("N135").Value = -16642.5006

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 16.923077
$ws.Range("I2").Value = 15.555555
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 15.555555
$ws.Range("L2").Value = 20
$ws.Range("M2").Value = 97.444445
$ws.Range("N2").Value = -246

$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10302

$ws.Range("H46").Value = 11166.667
$ws.Range("I46").Value = 5500
$ws.Range("J46").Value = 14000
$ws.Range("K46").Value = 5500
$ws.Range("L46").Value = 14000
$ws.Range("M46").Value = -5344
$ws.Range("N46").Value = -14312

$ws.Range("H57").Value = 19800
$ws.Range("J57").Value = 19800
$ws.Range("L57").Value = 19800
$ws.Range("N57").Value = -21440

$ws.Range("H132").Value = 254627.62
$ws.Range("I132").Value = 403804.6
$ws.Range("K132").Value = 1211413.8
$ws.Range("M132").Value = -1208883.8
